$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.439.65"
$ws.Range("E2").Value = "  -0.61%  "

$ws.Range("D3").Value = "1.955.33"
$ws.Range("E3").Value = "  -1.66%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.55%  "

$ws.Range("E6").Value = "  -1.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.84"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.39%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.366"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0857"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.57%  "

$ws.Range("E11").Value = "  +0.29%  "

$ws.Range("D12").Value = "2.241.20"
$ws.Range("E12").Value = "  -1.62%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -11.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.822"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.99%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.87%  "

$ws.Range("E16").Value = "  -5.14%  "

$ws.Range("D17").Value = "1.955.08"
$ws.Range("E17").Value = "  -1.53%  "

$ws.Range("D18").Value = "36.368.44"
$ws.Range("E18").Value = "  -0.53%  "

$ws.Range("D19").Value = "0.0₃0890"
$ws.Range("E19").Value = "  +2.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "230.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.22%  "

$ws.Range("E25").Value = "  -0.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.32%  "

$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.139"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.00%  "

$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.62%  "

$ws.Range("E30").Value = "  -1.88%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.95%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0653"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.73%  "

$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.51%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.16"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.82%  "

$ws.Range("E39").Value = "  -2.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0987"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.64%  "

$ws.Range("E41").Value = "  +0.49%  "

$ws.Range("E42").Value = "  -6.90%  "

$ws.Range("E43").Value = "  -1.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.30%  "

$ws.Range("D45").Value = "1.358.90"
$ws.Range("E45").Value = "  -1.28%  "

$ws.Range("E46").Value = "  -6.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.79%  "

$ws.Range("E49").Value = "  -0.93%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.29%  "

$ws.Range("D51").Value = "2.131.38"
$ws.Range("E51").Value = "  -1.94%  "
